$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C13").Value = "'02142"
Write-Output $ws1.Range("C13").Text
